# tag v 1.0.0 issue 2 issue 3
#
# The placeholder text in the "B.xlsx" template is renamed from the
# hyphenated form to the dotted form used by the templating engine:
#   {s-Name} -> {s.Name}
#   {s-Age}  -> {s.Age}
# "Merry Christmas!" is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "{s.Name}"
$ws.Range("E7").Value = "{s.Age}"

# Best-effort: the author's Excel session also moved the app window
# (bookViews/workbookView yWindow 1790 -> 2350) when the file was last
# saved. This is pure UI chrome with no effect on workbook content; set
# it anyway in case the host window object persists it.
$win = $excel.ActiveWindow
$win.Top = 2350
$win.Left = 240
